$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E14").Value = 8
$ws.Range("E15").Value = 5
$ws.Range("E16").Value = 8

$ws.Range("E15").Select()
